$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-21 Sunday" "2024-01-22 Monday"

Replace-Text "251÷6=" "646÷3="
Replace-Text "275÷2=" "554÷3="
Replace-Text "348÷9=" "618÷3="
Replace-Text "530÷6=" "539÷5="
Replace-Text "680÷3=" "935÷9="
Replace-Text "315÷9=" "640÷9="
Replace-Text "781÷8=" "281÷5="
Replace-Text "118÷9=" "318÷2="
Replace-Text "449÷3=" "654÷9="
Replace-Text "434÷5=" "289÷9="
Replace-Text "231÷3=" "134÷7="
Replace-Text "919÷2=" "613÷5="
Replace-Text "981÷2=" "638÷9="
Replace-Text "796÷8=" "221÷8="
Replace-Text "353÷7=" "598÷9="
Replace-Text "527÷6=" "788÷7="
Replace-Text "838÷2=" "882÷2="
Replace-Text "631÷3=" "587÷7="
Replace-Text "773÷3=" "808÷3="
Replace-Text "366÷8=" "238÷6="
Replace-Text "197÷4=" "761÷3="
Replace-Text "514÷2=" "699÷6="
Replace-Text "869÷9=" "512÷4="
Replace-Text "797÷4=" "683÷2="
Replace-Text "611÷3=" "972÷4="
